# Update "想去人数" (interested-count) values in column F across sheets,
# matching the output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 191
$ws1.Range("F3").Value = 404
$ws1.Range("F8").Value = 1065
$ws1.Range("F12").Value = 34
$ws1.Range("F13").Value = 312
$ws1.Range("F15").Value = 30
$ws1.Range("F16").Value = 61
$ws1.Range("F17").Value = 471
$ws1.Range("F18").Value = 445
$ws1.Range("F19").Value = 5597
$ws1.Range("F20").Value = 88
$ws1.Range("F21").Value = 1562
$ws1.Range("F22").Value = 368
$ws1.Range("F23").Value = 4771
$ws1.Range("F26").Value = 1499
$ws1.Range("F27").Value = 14
$ws1.Range("F28").Value = 24
$ws1.Range("F29").Value = 651
$ws1.Range("F30").Value = 63
$ws1.Range("F32").Value = 3792

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 10
$ws2.Range("F5").Value = 134
$ws2.Range("F8").Value = 101

# --- Sheet: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9383
$ws3.Range("F4").Value = 2128

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9383
$ws4.Range("F4").Value = 2128
$ws4.Range("F5").Value = 191
$ws4.Range("F6").Value = 404
$ws4.Range("F11").Value = 1065
$ws4.Range("F14").Value = 34
$ws4.Range("F15").Value = 312
$ws4.Range("F17").Value = 30
$ws4.Range("F18").Value = 61
$ws4.Range("F22").Value = 445
$ws4.Range("F23").Value = 5597
$ws4.Range("F24").Value = 88
$ws4.Range("F25").Value = 1562
$ws4.Range("F28").Value = 368
$ws4.Range("F31").Value = 4771
$ws4.Range("F34").Value = 1499
$ws4.Range("F35").Value = 14
$ws4.Range("F36").Value = 24
$ws4.Range("F37").Value = 651
$ws4.Range("F38").Value = 63
$ws4.Range("F46").Value = 3792
